# Auto-generated Excel COM-interop script
# Applies "Add data for 2024-08-23" updates to violent-crime-full-year.xlsx
# For each affected worksheet, updates the 2024 (column K) values (and a couple of
# column J corrections) plus adds one newly-populated cell (Hyde Park K5).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 5149
$ws.Range("K3").Value = 5288
$ws.Range("K4").Value = 1104
$ws.Range("K5").Value = 379
$ws.Range("J6").Value = 11053
$ws.Range("K6").Value = 5927
$ws.Range("J7").Value = 29297
$ws.Range("K7").Value = 17847

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 153
$ws.Range("K4").Value = 68
$ws.Range("K5").Value = 41
$ws.Range("K6").Value = 132
$ws.Range("K7").Value = 533
$ws.Range("K8").Value = 1210
$ws.Range("K11").Value = 343
$ws.Range("K14").Value = 94
$ws.Range("K15").Value = 180
$ws.Range("K17").Value = 35
$ws.Range("K18").Value = 121
$ws.Range("K19").Value = 529
$ws.Range("K20").Value = 410
$ws.Range("K24").Value = 53
$ws.Range("K25").Value = 86
$ws.Range("K27").Value = 164
$ws.Range("K29").Value = 957
$ws.Range("K33").Value = 757
$ws.Range("K42").Value = 657
$ws.Range("K43").Value = 155
$ws.Range("K44").Value = 156
$ws.Range("K48").Value = 223
$ws.Range("K49").Value = 101
$ws.Range("K50").Value = 86
$ws.Range("J52").Value = 743
$ws.Range("K52").Value = 466
$ws.Range("K53").Value = 229
$ws.Range("K54").Value = 351
$ws.Range("K63").Value = 54
$ws.Range("K64").Value = 112
$ws.Range("K67").Value = 680
$ws.Range("K68").Value = 47
$ws.Range("K75").Value = 62
$ws.Range("K79").Value = 441
$ws.Range("K83").Value = 397
$ws.Range("K85").Value = 838
$ws.Range("K86").Value = 120
$ws.Range("K88").Value = 198
$ws.Range("K89").Value = 261
$ws.Range("K92").Value = 67
$ws.Range("K99").Value = 303
$ws.Range("J101").Value = 29297
$ws.Range("K101").Value = 17847

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 94

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 182
$ws.Range("K7").Value = 533

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 114
$ws.Range("K6").Value = 117
$ws.Range("K7").Value = 343

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 73
$ws.Range("K6").Value = 80
$ws.Range("K7").Value = 261

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 283
$ws.Range("K3").Value = 283
$ws.Range("K5").Value = 25
$ws.Range("K7").Value = 838

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K4").Value = 26
$ws.Range("J6").Value = 316
$ws.Range("J7").Value = 743
$ws.Range("K7").Value = 466

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 58
$ws.Range("K3").Value = 56
$ws.Range("K7").Value = 229

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K3").Value = 362
$ws.Range("K4").Value = 70
$ws.Range("K6").Value = 411
$ws.Range("K7").Value = 1210

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 135
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 397

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 204
$ws.Range("K3").Value = 282
$ws.Range("K6").Value = 220
$ws.Range("K7").Value = 757

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K6").Value = 77
$ws.Range("K7").Value = 303

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 198
$ws.Range("K3").Value = 239
$ws.Range("K7").Value = 680

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 90
$ws.Range("K6").Value = 186
$ws.Range("K7").Value = 351

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 276
$ws.Range("K3").Value = 341
$ws.Range("K6").Value = 265
$ws.Range("K7").Value = 957

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K4").Value = 32
$ws.Range("K6").Value = 111
$ws.Range("K7").Value = 223

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 156
$ws.Range("K6").Value = 168
$ws.Range("K7").Value = 529

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K3").Value = 42
$ws.Range("K7").Value = 156

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K3").Value = 36
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 174
$ws.Range("K3").Value = 203
$ws.Range("K6").Value = 249
$ws.Range("K7").Value = 657

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K3").Value = 16
$ws.Range("K4").Value = 7
$ws.Range("K6").Value = 48

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 146
$ws.Range("K3").Value = 141
$ws.Range("K4").Value = 29
$ws.Range("K7").Value = 441

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 112

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 132
$ws.Range("K4").Value = 19
$ws.Range("K6").Value = 118
$ws.Range("K7").Value = 410

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 33
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("K2").Value = 13
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 63
$ws.Range("K7").Value = 180

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 50
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K6").Value = 80
$ws.Range("K7").Value = 198

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 44
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 164

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K3").Value = 21
$ws.Range("K4").Value = 51
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K5").Value = 1
$ws.Range("K7").Value = 155

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 68
